$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Paragraph 4 ("A brief explanation ...") - expand the sentence and add
# a second sentence about documenting design decisions.
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "A brief explanation of how you went about measuring the latency and throughput.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A brief but detailed explanation of how you went about measuring the latency and throughput. You should include any decisions you made (e.g., wall clock vs. CPU clock, synchronous vs. asynchronous RPCs, etc.).",
    2) | Out-Null

# -----------------------------------------------------------------------
# Paragraph 5 ("Paste the C++ code ...") - replace entirely with the new
# "compress your source code" instructions.
# -----------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Paste the C++ code*") {
        $target = $d.Paragraphs.Item($i)
    }
}

$pStart = $target.Range.Start
$pEnd = $target.Range.End
# Keep the trailing paragraph mark untouched; only replace the visible text.
$textEnd = $pEnd - 1
$bodyRng = $d.Range($pStart, $textEnd)
$bodyRng.Text = "Compress all your relevant source code (only the files you actually created/modified; not the whole gRPC repository!) along with the PDF in a zip archive or tarball and submit to CatCourses. "

# -----------------------------------------------------------------------
# New sub-bullet paragraph right after it, one list level deeper. Insert
# this (plain, unbolded) paragraph before bolding "relevant" so it never
# inherits that character formatting.
# -----------------------------------------------------------------------
$p5 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Compress all your relevant source code*") {
        $p5 = $d.Paragraphs.Item($i)
    }
}

$insertPt = $d.Range($p5.Range.End, $p5.Range.End)
$insertPt.InsertParagraphAfter()
$p6 = $p5.Next()

$p6.Range.Text = "The files you modify for this lab will probably be greeter_client.cpp and greeter_server.cpp, though it depends."
$p6.Range.ListFormat.ListLevelNumber = 2

# Now bold just the word "relevant" in paragraph 5 (done last so the
# bold character formatting cannot leak into the new paragraph below).
$boldStart = $pStart + "Compress all your ".Length
$boldEnd = $boldStart + "relevant".Length
$boldRng = $d.Range($boldStart, $boldEnd)
$boldRng.Bold = 1
